$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.26275633660550568
$ws.Range("A2").Value = -0.0099999997061956947
$ws.Range("A3").Value = -0.0089999997043062052
$ws.Range("A4").Value = 0.061996661545101972
$ws.Range("A5").Value = -0.0059999997109549952
$ws.Range("A6").Value = -0.0059999997030324437
$ws.Range("A7").Value = -0.019999999648014466
$ws.Range("A8").Value = -0.019999999648648625
$ws.Range("A9").Value = -0.0059999997060842247
$ws.Range("A10").Value = -0.0059999997089974499
$ws.Range("A11").Value = -0.0044999997152856963
$ws.Range("A12").Value = 0.065372065473721097
$ws.Range("A13").Value = -0.0059999997096626956
$ws.Range("A14").Value = -0.011999999685946783
$ws.Range("A15").Value = 0.029047236886433403
$ws.Range("A16").Value = -0.0059999997088362456
$ws.Range("A17").Value = -0.0059999997076705114
$ws.Range("A18").Value = -0.0089999996956278139
$ws.Range("A19").Value = -0.0089999997095042694
$ws.Range("A20").Value = -0.046721325296335436
$ws.Range("A21").Value = -0.0089999997025191902
$ws.Range("A22").Value = -0.0089999997022234268
$ws.Range("A23").Value = -0.0089999997016425581
$ws.Range("A24").Value = -0.041999999567297941
$ws.Range("A25").Value = -0.041999999564986013
$ws.Range("A26").Value = -0.0059999997015083295
$ws.Range("A27").Value = -0.029392592521128247
$ws.Range("A28").Value = -0.0059999996925714782
$ws.Range("A29").Value = -0.01199999966422638
$ws.Range("A30").Value = -0.019999999630496479
$ws.Range("A31").Value = -0.01499999964702603
$ws.Range("A32").Value = -0.020999999622959287
$ws.Range("A33").Value = -0.0059999996817685641
